$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "testting"
$ws.Cells.Item(2, 2).Value = "2024-10-01 22:09:12"
